$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as plain text, even if the text looks like a
# date/number, without leaving behind any new/duplicated cell style.
# We do this by temporarily writing a text-returning formula (T() always
# returns text) and then converting the cell to a static value with
# PasteSpecial (values only), which keeps the existing cell style intact.
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = "=T(""$escaped"")"
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $range.Application.CutCopyMode = $false
}

# Helper: set a brand-new cell's value + copy the number/cell formatting
# (style) from a reference cell that already has the desired style, without
# touching the reference cell's own content.
function Set-NewCell($range, [string]$text, $formatSource) {
    $range.Value = $text
    $formatSource.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
    $range.Application.CutCopyMode = $false
}

# 1. Ativação / activation date: 01/01/2012 -> 01/01/2023 (rows 8 and 15)
Set-TextValue $ws.Range("B8") "01/01/2023"
Set-TextValue $ws.Range("C8") "01/01/2023"
Set-TextValue $ws.Range("B15") "01/01/2023"
Set-TextValue $ws.Range("C15") "01/01/2023"

# 2. Objectives text (new cells B11 / C11), formatted like B13 / C13
Set-NewCell $ws.Range("B11") "Provide the student with knowledge of the main techniques of physical and chemical characterization of materials." $ws.Range("B13")
Set-NewCell $ws.Range("C11") "Provide the student with knowledge of the main techniques of physical and chemical characterization of materials." $ws.Range("C13")

# 3. Short syllabus text (new cells B14 / C14)
Set-NewCell $ws.Range("B14") "Granulometric and surface analysis. Microstructural analyses. Thermal analysis. Rheometry." $ws.Range("B13")
Set-NewCell $ws.Range("C14") "Granulometric and surface analysis. Microstructural analyses. Thermal analysis. Rheometry." $ws.Range("C13")

# 4. Full syllabus text (new cells B16 / C16)
Set-NewCell $ws.Range("B16") "Grain size analysis. BET adsorption, porosity and pycnometry.Microstructural analysis: X-ray diffraction, Laue figure; X-ray scattering (SAXS). Electron diffraction. Optical Microscopy. Electron microscopy, X-ray microanalysis (EDX and WDX).Thermal analysis: Differential thermal analysis (DTA), differential scanning calorimetry (DSC) and thermogravimetry (TGA).Rheometry of liquids, solutions and pastes." $ws.Range("B13")
Set-NewCell $ws.Range("C16") "Grain size analysis. BET adsorption, porosity and pycnometry.Microstructural analysis: X-ray diffraction, Laue figure; X-ray scattering (SAXS). Electron diffraction. Optical Microscopy. Electron microscopy, X-ray microanalysis (EDX and WDX).Thermal analysis: Differential thermal analysis (DTA), differential scanning calorimetry (DSC) and thermogravimetry (TGA).Rheometry of liquids, solutions and pastes." $ws.Range("C13")
